$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers I1 ("I0") and J1 ("IF"), matching the style (bold, bordered,
# centered) already used by the other header cells in row 1.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New data columns I ("I0") and J ("IF") for rows 2-38.
$data = @(
    @(4, 6),
    @(7, 9),
    @(7, 8),
    @(1, 6),
    @(1, 3),
    @(1, 4),
    @(1, 3),
    @(1, 5),
    @(1, 6),
    @(1, 5),
    @(1, 7),
    @(1, 7),
    @(1, 6),
    @(1, 5),
    @(1, 5),
    @(1, 6),
    @(1, 6),
    @(1, 6),
    @(1, 5),
    @(1, 3),
    @(1, 9),
    @(1, 7),
    @(1, 7),
    @(1, 5),
    @(1, 8),
    @(1, 7),
    @(1, 6),
    @(1, 6),
    @(1, 6),
    @(1, 3),
    @(1, 6),
    @(1, 7),
    @(6, 9),
    @(1, 5),
    @(1, 4),
    @(1, 3),
    @(1, 2)
)

for ($idx = 0; $idx -lt $data.Length; $idx++) {
    $row = $idx + 2
    $pair = $data[$idx]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}

